# The workbook stores "ObjTables" header metadata as plain text in the
# first cell(s) of every sheet (these sheets are all named "!!<TableId>").
# This commit adds a schema='SBtab' attribute to every such header string,
# reorders a couple of attributes, and bumps the embedded timestamp.
#
# Sheet "!!Compartment" is special: its A1 holds the *document*-level
# "!!!ObjTables ..." header and its A2 holds the *table*-level
# "!!ObjTables ... id='Compartment' ..." header. Every other sheet only
# has the table-level header, in A1.

$wb = $excel.ActiveWorkbook

$docHeaderOld = "!!!ObjTables objTablesVersion='0.0.8' date='2020-03-09 15:31:26'"
$docHeaderNew = "!!!ObjTables schema='SBtab' objTablesVersion='0.0.8' date='2020-03-09 23:58:48'"

# One entry per data table / worksheet, in workbook (tab) order.
$tables = @(
    @{ Name = "Compartment";            Date = "2020-03-09 23:58:48" },
    @{ Name = "Compound";               Date = "2020-03-09 23:58:48" },
    @{ Name = "Definition";             Date = "2020-03-09 23:58:48" },
    @{ Name = "Enzyme";                 Date = "2020-03-09 23:58:48" },
    @{ Name = "FbcObjective";           Date = "2020-03-09 23:58:48" },
    @{ Name = "Gene";                   Date = "2020-03-09 23:58:48" },
    @{ Name = "Layout";                 Date = "2020-03-09 23:58:48" },
    @{ Name = "Measurement";            Date = "2020-03-09 23:58:48" },
    @{ Name = "PbConfig";               Date = "2020-03-09 23:58:48" },
    @{ Name = "Position";               Date = "2020-03-09 23:58:48" },
    @{ Name = "Protein";                Date = "2020-03-09 23:58:48" },
    @{ Name = "Quantity";               Date = "2020-03-09 23:58:48" },
    @{ Name = "QuantityInfo";           Date = "2020-03-09 23:58:48" },
    @{ Name = "QuantityMatrix";         Date = "2020-03-09 23:58:48" },
    @{ Name = "Reaction";               Date = "2020-03-09 23:58:48" },
    @{ Name = "ReactionStoichiometry";  Date = "2020-03-09 23:58:48" },
    @{ Name = "Regulator";              Date = "2020-03-09 23:58:49" },
    @{ Name = "Relation";               Date = "2020-03-09 23:58:49" },
    @{ Name = "Relationship";           Date = "2020-03-09 23:58:49" },
    @{ Name = "SparseMatrix";           Date = "2020-03-09 23:58:49" },
    @{ Name = "SparseMatrixColumn";     Date = "2020-03-09 23:58:49" },
    @{ Name = "SparseMatrixOrdered";    Date = "2020-03-09 23:58:49" },
    @{ Name = "SparseMatrixRow";        Date = "2020-03-09 23:58:49" },
    @{ Name = "StoichiometricMatrix";   Date = "2020-03-09 23:58:49" },
    @{ Name = "rxnconContingencyList";  Date = "2020-03-09 23:58:49" },
    @{ Name = "rxnconReactionList";     Date = "2020-03-09 23:58:49" }
)

foreach ($t in $tables) {
    $ws = $wb.Worksheets.Item("!!" + $t.Name)

    # Sheets are protected (read-only); temporarily unprotect to write.
    $ws.Unprotect()

    $tableHeaderNew = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='" + $t.Name + "' name='" + $t.Name + "' date='" + $t.Date + "' objTablesVersion='0.0.8'"

    if ($t.Name -eq "Compartment") {
        # This sheet carries both the document-level header (row 1) and
        # the table-level header (row 2).
        $ws.Range("A1").Value = $docHeaderNew
        $ws.Range("A2").Value = $tableHeaderNew
    } else {
        $ws.Range("A1").Value = $tableHeaderNew
    }

    $ws.Protect()
}
